$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Countries re-ranked by case count (some rows shifted), and the
# Casos totales / Nuevos casos / Casos activos / Recuperados / Casos
# criticos / Muertes hoy / Muertes columns (B:H) refreshed with new data.

$ws.Cells.Item(5, 1).Value = 'Italia'
$ws.Cells.Item(5, 2).Value = 53578
$ws.Cells.Item(5, 3).Value = 6557
$ws.Cells.Item(5, 4).Value = 6072
$ws.Cells.Item(5, 5).Value = 42681
$ws.Cells.Item(5, 6).Value = 2857
$ws.Cells.Item(5, 7).Value = 793
$ws.Cells.Item(5, 8).Value = 4825

$ws.Cells.Item(6, 1).Value = 'España'
$ws.Cells.Item(6, 2).Value = 25374
$ws.Cells.Item(6, 3).Value = 3803
$ws.Cells.Item(6, 4).Value = 2125
$ws.Cells.Item(6, 5).Value = 21871
$ws.Cells.Item(6, 6).Value = 1612
$ws.Cells.Item(6, 7).Value = 285
$ws.Cells.Item(6, 8).Value = 1378

$ws.Cells.Item(7, 1).Value = 'Estados Unidos'
$ws.Cells.Item(7, 2).Value = 22132
$ws.Cells.Item(7, 3).Value = 2749
$ws.Cells.Item(7, 4).Value = 147
$ws.Cells.Item(7, 5).Value = 21703
$ws.Cells.Item(7, 6).Value = 64
$ws.Cells.Item(7, 7).Value = 26
$ws.Cells.Item(7, 8).Value = 282

$ws.Cells.Item(8, 1).Value = 'Alemania'
$ws.Cells.Item(8, 2).Value = 21854
$ws.Cells.Item(8, 3).Value = 2006
$ws.Cells.Item(8, 4).Value = 209
$ws.Cells.Item(8, 5).Value = 21568
$ws.Cells.Item(8, 6).Value = 2
$ws.Cells.Item(8, 7).Value = 9
$ws.Cells.Item(8, 8).Value = 77

$ws.Cells.Item(12, 1).Value = 'Suiza'
$ws.Cells.Item(12, 2).Value = 6371
$ws.Cells.Item(12, 3).Value = 756
$ws.Cells.Item(12, 4).Value = 15
$ws.Cells.Item(12, 5).Value = 6284
$ws.Cells.Item(12, 6).Value = 141
$ws.Cells.Item(12, 7).Value = 16
$ws.Cells.Item(12, 8).Value = 72

$ws.Cells.Item(17, 1).Value = 'Noruega'
$ws.Cells.Item(17, 2).Value = 2093
$ws.Cells.Item(17, 3).Value = 134
$ws.Cells.Item(17, 4).Value = 1
$ws.Cells.Item(17, 5).Value = 2085
$ws.Cells.Item(17, 6).Value = 28
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = 7

$ws.Cells.Item(62, 1).Value = 'Serbia'
$ws.Cells.Item(62, 2).Value = 171
$ws.Cells.Item(62, 3).Value = 36
$ws.Cells.Item(62, 4).Value = 2
$ws.Cells.Item(62, 5).Value = 168
$ws.Cells.Item(62, 6).Value = 4
$ws.Cells.Item(62, 7).Value = 0
$ws.Cells.Item(62, 8).Value = 1

$ws.Cells.Item(63, 1).Value = 'Bulgaria'
$ws.Cells.Item(63, 2).Value = 163
$ws.Cells.Item(63, 3).Value = 36
$ws.Cells.Item(63, 4).Value = 3
$ws.Cells.Item(63, 5).Value = 157
$ws.Cells.Item(63, 6).Value = 3
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 3

$ws.Cells.Item(64, 1).Value = 'Armenia'
$ws.Cells.Item(64, 2).Value = 160
$ws.Cells.Item(64, 3).Value = 24
$ws.Cells.Item(64, 4).Value = 1
$ws.Cells.Item(64, 5).Value = 159
$ws.Cells.Item(64, 6).Value = 2
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 8).Value = 0

$ws.Cells.Item(65, 1).Value = 'Argentina'
$ws.Cells.Item(65, 2).Value = 158
$ws.Cells.Item(65, 3).Value = 0
$ws.Cells.Item(65, 4).Value = 3
$ws.Cells.Item(65, 5).Value = 151
$ws.Cells.Item(65, 6).Value = 0
$ws.Cells.Item(65, 7).Value = 1
$ws.Cells.Item(65, 8).Value = 4

$ws.Cells.Item(66, 1).Value = 'Taiwan'
$ws.Cells.Item(66, 2).Value = 153
$ws.Cells.Item(66, 3).Value = 18
$ws.Cells.Item(66, 4).Value = 28
$ws.Cells.Item(66, 5).Value = 123
$ws.Cells.Item(66, 6).Value = 0
$ws.Cells.Item(66, 7).Value = 0
$ws.Cells.Item(66, 8).Value = 2

$ws.Cells.Item(67, 1).Value = 'Emiratos Arabes Unidos'
$ws.Cells.Item(67, 2).Value = 153
$ws.Cells.Item(67, 3).Value = 13
$ws.Cells.Item(67, 4).Value = 38
$ws.Cells.Item(67, 5).Value = 113
$ws.Cells.Item(67, 6).Value = 2
$ws.Cells.Item(67, 7).Value = 0
$ws.Cells.Item(67, 8).Value = 2

$ws.Cells.Item(68, 1).Value = 'San Marino'
$ws.Cells.Item(68, 2).Value = 151
$ws.Cells.Item(68, 3).Value = 0
$ws.Cells.Item(68, 4).Value = 4
$ws.Cells.Item(68, 5).Value = 127
$ws.Cells.Item(68, 6).Value = 12
$ws.Cells.Item(68, 7).Value = 6
$ws.Cells.Item(68, 8).Value = 20

$ws.Cells.Item(81, 1).Value = 'Republica de Macedonia'
$ws.Cells.Item(81, 2).Value = 85
$ws.Cells.Item(81, 3).Value = 9
$ws.Cells.Item(81, 4).Value = 1
$ws.Cells.Item(81, 5).Value = 84
$ws.Cells.Item(81, 6).Value = 1
$ws.Cells.Item(81, 7).Value = 0
$ws.Cells.Item(81, 8).Value = 0

$ws.Cells.Item(82, 1).Value = 'Jordania'
$ws.Cells.Item(82, 2).Value = 85
$ws.Cells.Item(82, 3).Value = 0
$ws.Cells.Item(82, 4).Value = 1
$ws.Cells.Item(82, 5).Value = 84
$ws.Cells.Item(82, 6).Value = 0
$ws.Cells.Item(82, 7).Value = 0
$ws.Cells.Item(82, 8).Value = 0

$ws.Cells.Item(116, 1).Value = 'Ghana'
$ws.Cells.Item(116, 2).Value = 19
$ws.Cells.Item(116, 3).Value = 3
$ws.Cells.Item(116, 4).Value = 0
$ws.Cells.Item(116, 5).Value = 19
$ws.Cells.Item(116, 6).Value = 0
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 0

$ws.Cells.Item(117, 1).Value = 'Bolivia'
$ws.Cells.Item(117, 2).Value = 19
$ws.Cells.Item(117, 3).Value = 3
$ws.Cells.Item(117, 4).Value = 0
$ws.Cells.Item(117, 5).Value = 19
$ws.Cells.Item(117, 6).Value = 0
$ws.Cells.Item(117, 7).Value = 0
$ws.Cells.Item(117, 8).Value = 0

$ws.Cells.Item(121, 1).Value = 'Ruanda'
$ws.Cells.Item(121, 2).Value = 17
$ws.Cells.Item(121, 3).Value = 0
$ws.Cells.Item(121, 4).Value = 0
$ws.Cells.Item(121, 5).Value = 17
$ws.Cells.Item(121, 6).Value = 0
$ws.Cells.Item(121, 7).Value = 0
$ws.Cells.Item(121, 8).Value = 0

$ws.Cells.Item(122, 1).Value = 'Puerto Rico'
$ws.Cells.Item(122, 2).Value = 17
$ws.Cells.Item(122, 3).Value = 3
$ws.Cells.Item(122, 4).Value = 0
$ws.Cells.Item(122, 5).Value = 17
$ws.Cells.Item(122, 6).Value = 0
$ws.Cells.Item(122, 7).Value = 0
$ws.Cells.Item(122, 8).Value = 0

$ws.Cells.Item(125, 1).Value = 'Guayana Francesa'
$ws.Cells.Item(125, 2).Value = 15
$ws.Cells.Item(125, 3).Value = 0
$ws.Cells.Item(125, 4).Value = 0
$ws.Cells.Item(125, 5).Value = 15
$ws.Cells.Item(125, 6).Value = 0
$ws.Cells.Item(125, 7).Value = 0
$ws.Cells.Item(125, 8).Value = 0

$ws.Cells.Item(126, 1).Value = 'Polinesia Francesa'
$ws.Cells.Item(126, 2).Value = 15
$ws.Cells.Item(126, 3).Value = 4
$ws.Cells.Item(126, 4).Value = 0
$ws.Cells.Item(126, 5).Value = 15
$ws.Cells.Item(126, 6).Value = 0
$ws.Cells.Item(126, 7).Value = 0
$ws.Cells.Item(126, 8).Value = 0

$ws.Cells.Item(127, 1).Value = 'Kirguistan'
$ws.Cells.Item(127, 2).Value = 14
$ws.Cells.Item(127, 3).Value = 8
$ws.Cells.Item(127, 4).Value = 0
$ws.Cells.Item(127, 5).Value = 14
$ws.Cells.Item(127, 6).Value = 0
$ws.Cells.Item(127, 7).Value = 0
$ws.Cells.Item(127, 8).Value = 0

$ws.Cells.Item(128, 1).Value = 'Montenegro'
$ws.Cells.Item(128, 2).Value = 14
$ws.Cells.Item(128, 3).Value = 0
$ws.Cells.Item(128, 4).Value = 0
$ws.Cells.Item(128, 5).Value = 14
$ws.Cells.Item(128, 6).Value = 0
$ws.Cells.Item(128, 7).Value = 0
$ws.Cells.Item(128, 8).Value = 0

$ws.Cells.Item(129, 1).Value = 'Costa de Marfil'
$ws.Cells.Item(129, 2).Value = 14
$ws.Cells.Item(129, 3).Value = 5
$ws.Cells.Item(129, 4).Value = 1
$ws.Cells.Item(129, 5).Value = 13
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 0

$ws.Cells.Item(130, 1).Value = 'Mauricio'
$ws.Cells.Item(130, 2).Value = 14
$ws.Cells.Item(130, 3).Value = 2
$ws.Cells.Item(130, 4).Value = 0
$ws.Cells.Item(130, 5).Value = 13
$ws.Cells.Item(130, 6).Value = 0
$ws.Cells.Item(130, 7).Value = 1
$ws.Cells.Item(130, 8).Value = 1

$ws.Cells.Item(137, 1).Value = 'Seychelles'
$ws.Cells.Item(137, 2).Value = 7
$ws.Cells.Item(137, 3).Value = 0
$ws.Cells.Item(137, 4).Value = 0
$ws.Cells.Item(137, 5).Value = 7
$ws.Cells.Item(137, 6).Value = 0
$ws.Cells.Item(137, 7).Value = 0
$ws.Cells.Item(137, 8).Value = 0

$ws.Cells.Item(138, 1).Value = 'Kenia'
$ws.Cells.Item(138, 2).Value = 7
$ws.Cells.Item(138, 3).Value = 0
$ws.Cells.Item(138, 4).Value = 0
$ws.Cells.Item(138, 5).Value = 7
$ws.Cells.Item(138, 6).Value = 0
$ws.Cells.Item(138, 7).Value = 0
$ws.Cells.Item(138, 8).Value = 0

$ws.Cells.Item(139, 1).Value = 'Mayotte'
$ws.Cells.Item(139, 2).Value = 7
$ws.Cells.Item(139, 3).Value = 0
$ws.Cells.Item(139, 4).Value = 0
$ws.Cells.Item(139, 5).Value = 7
$ws.Cells.Item(139, 6).Value = 0
$ws.Cells.Item(139, 7).Value = 0
$ws.Cells.Item(139, 8).Value = 0

$ws.Cells.Item(140, 1).Value = 'Islas Virgenes de los Estados Unidos'
$ws.Cells.Item(140, 2).Value = 6
$ws.Cells.Item(140, 3).Value = 3
$ws.Cells.Item(140, 4).Value = 0
$ws.Cells.Item(140, 5).Value = 6
$ws.Cells.Item(140, 6).Value = 0
$ws.Cells.Item(140, 7).Value = 0
$ws.Cells.Item(140, 8).Value = 0

$ws.Cells.Item(141, 1).Value = 'Barbados'
$ws.Cells.Item(141, 2).Value = 6
$ws.Cells.Item(141, 3).Value = 0
$ws.Cells.Item(141, 4).Value = 0
$ws.Cells.Item(141, 5).Value = 6
$ws.Cells.Item(141, 6).Value = 0
$ws.Cells.Item(141, 7).Value = 0
$ws.Cells.Item(141, 8).Value = 0

$ws.Cells.Item(142, 1).Value = 'Guinea Ecuatorial'
$ws.Cells.Item(142, 2).Value = 6
$ws.Cells.Item(142, 3).Value = 0
$ws.Cells.Item(142, 4).Value = 0
$ws.Cells.Item(142, 5).Value = 6
$ws.Cells.Item(142, 6).Value = 0
$ws.Cells.Item(142, 7).Value = 0
$ws.Cells.Item(142, 8).Value = 0

$ws.Cells.Item(143, 1).Value = 'Tanzania'
$ws.Cells.Item(143, 2).Value = 6
$ws.Cells.Item(143, 3).Value = 0
$ws.Cells.Item(143, 4).Value = 0
$ws.Cells.Item(143, 5).Value = 6
$ws.Cells.Item(143, 6).Value = 0
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = 0

$ws.Cells.Item(144, 1).Value = 'Guyana'
$ws.Cells.Item(144, 2).Value = 5
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(144, 4).Value = 0
$ws.Cells.Item(144, 5).Value = 4
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 1

$ws.Cells.Item(145, 1).Value = 'Aruba'
$ws.Cells.Item(145, 2).Value = 5
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(145, 4).Value = 1
$ws.Cells.Item(145, 5).Value = 4
$ws.Cells.Item(145, 6).Value = 0
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 0

$ws.Cells.Item(151, 1).Value = 'Congo'
$ws.Cells.Item(151, 2).Value = 3
$ws.Cells.Item(151, 3).Value = 0
$ws.Cells.Item(151, 4).Value = 0
$ws.Cells.Item(151, 5).Value = 3
$ws.Cells.Item(151, 6).Value = 0
$ws.Cells.Item(151, 7).Value = 0
$ws.Cells.Item(151, 8).Value = 0

$ws.Cells.Item(152, 1).Value = 'Liberia'
$ws.Cells.Item(152, 2).Value = 3
$ws.Cells.Item(152, 3).Value = 1
$ws.Cells.Item(152, 4).Value = 0
$ws.Cells.Item(152, 5).Value = 3
$ws.Cells.Item(152, 6).Value = 0
$ws.Cells.Item(152, 7).Value = 0
$ws.Cells.Item(152, 8).Value = 0

$ws.Cells.Item(153, 1).Value = 'Republica de Africa Central'
$ws.Cells.Item(153, 2).Value = 3
$ws.Cells.Item(153, 3).Value = 0
$ws.Cells.Item(153, 4).Value = 0
$ws.Cells.Item(153, 5).Value = 3
$ws.Cells.Item(153, 6).Value = 0
$ws.Cells.Item(153, 7).Value = 0
$ws.Cells.Item(153, 8).Value = 0

$ws.Cells.Item(154, 1).Value = 'Zimbabue'
$ws.Cells.Item(154, 2).Value = 3
$ws.Cells.Item(154, 3).Value = 2
$ws.Cells.Item(154, 4).Value = 0
$ws.Cells.Item(154, 5).Value = 3
$ws.Cells.Item(154, 6).Value = 0
$ws.Cells.Item(154, 7).Value = 0
$ws.Cells.Item(154, 8).Value = 0

$ws.Cells.Item(155, 1).Value = 'San Bartolome'
$ws.Cells.Item(155, 2).Value = 3
$ws.Cells.Item(155, 3).Value = 0
$ws.Cells.Item(155, 4).Value = 0
$ws.Cells.Item(155, 5).Value = 3
$ws.Cells.Item(155, 6).Value = 0
$ws.Cells.Item(155, 7).Value = 0
$ws.Cells.Item(155, 8).Value = 0

$ws.Cells.Item(156, 1).Value = 'Madagascar'
$ws.Cells.Item(156, 2).Value = 3
$ws.Cells.Item(156, 3).Value = 0
$ws.Cells.Item(156, 4).Value = 0
$ws.Cells.Item(156, 5).Value = 3
$ws.Cells.Item(156, 6).Value = 0
$ws.Cells.Item(156, 7).Value = 0
$ws.Cells.Item(156, 8).Value = 0

$ws.Cells.Item(157, 1).Value = 'El Salvador'
$ws.Cells.Item(157, 2).Value = 3
$ws.Cells.Item(157, 3).Value = 2
$ws.Cells.Item(157, 4).Value = 0
$ws.Cells.Item(157, 5).Value = 3
$ws.Cells.Item(157, 6).Value = 0
$ws.Cells.Item(157, 7).Value = 0
$ws.Cells.Item(157, 8).Value = 0

$ws.Cells.Item(158, 1).Value = 'Namibia'
$ws.Cells.Item(158, 2).Value = 3
$ws.Cells.Item(158, 3).Value = 0
$ws.Cells.Item(158, 4).Value = 0
$ws.Cells.Item(158, 5).Value = 3
$ws.Cells.Item(158, 6).Value = 0
$ws.Cells.Item(158, 7).Value = 0
$ws.Cells.Item(158, 8).Value = 0

$ws.Cells.Item(159, 1).Value = 'Cabo Verde'
$ws.Cells.Item(159, 2).Value = 3
$ws.Cells.Item(159, 3).Value = 2
$ws.Cells.Item(159, 4).Value = 0
$ws.Cells.Item(159, 5).Value = 3
$ws.Cells.Item(159, 6).Value = 0
$ws.Cells.Item(159, 7).Value = 0
$ws.Cells.Item(159, 8).Value = 0

$ws.Cells.Item(162, 1).Value = 'Groenlandia'
$ws.Cells.Item(162, 2).Value = 2
$ws.Cells.Item(162, 3).Value = 0
$ws.Cells.Item(162, 4).Value = 0
$ws.Cells.Item(162, 5).Value = 2
$ws.Cells.Item(162, 6).Value = 0
$ws.Cells.Item(162, 7).Value = 0
$ws.Cells.Item(162, 8).Value = 0

$ws.Cells.Item(163, 1).Value = 'Benin'
$ws.Cells.Item(163, 2).Value = 2
$ws.Cells.Item(163, 3).Value = 0
$ws.Cells.Item(163, 4).Value = 0
$ws.Cells.Item(163, 5).Value = 2
$ws.Cells.Item(163, 6).Value = 0
$ws.Cells.Item(163, 7).Value = 0
$ws.Cells.Item(163, 8).Value = 0

$ws.Cells.Item(164, 1).Value = 'Bermudas'
$ws.Cells.Item(164, 2).Value = 2
$ws.Cells.Item(164, 3).Value = 0
$ws.Cells.Item(164, 4).Value = 0
$ws.Cells.Item(164, 5).Value = 2
$ws.Cells.Item(164, 6).Value = 0
$ws.Cells.Item(164, 7).Value = 0
$ws.Cells.Item(164, 8).Value = 0

$ws.Cells.Item(165, 1).Value = 'Mauritania'
$ws.Cells.Item(165, 2).Value = 2
$ws.Cells.Item(165, 3).Value = 0
$ws.Cells.Item(165, 4).Value = 0
$ws.Cells.Item(165, 5).Value = 2
$ws.Cells.Item(165, 6).Value = 0
$ws.Cells.Item(165, 7).Value = 0
$ws.Cells.Item(165, 8).Value = 0

$ws.Cells.Item(166, 1).Value = 'Guinea'
$ws.Cells.Item(166, 2).Value = 2
$ws.Cells.Item(166, 3).Value = 0
$ws.Cells.Item(166, 4).Value = 0
$ws.Cells.Item(166, 5).Value = 2
$ws.Cells.Item(166, 6).Value = 0
$ws.Cells.Item(166, 7).Value = 0
$ws.Cells.Item(166, 8).Value = 0

$ws.Cells.Item(167, 1).Value = 'Isla de Man'
$ws.Cells.Item(167, 2).Value = 2
$ws.Cells.Item(167, 3).Value = 0
$ws.Cells.Item(167, 4).Value = 0
$ws.Cells.Item(167, 5).Value = 2
$ws.Cells.Item(167, 6).Value = 0
$ws.Cells.Item(167, 7).Value = 0
$ws.Cells.Item(167, 8).Value = 0

$ws.Cells.Item(168, 1).Value = 'Haiti'
$ws.Cells.Item(168, 2).Value = 2
$ws.Cells.Item(168, 3).Value = 0
$ws.Cells.Item(168, 4).Value = 0
$ws.Cells.Item(168, 5).Value = 2
$ws.Cells.Item(168, 6).Value = 0
$ws.Cells.Item(168, 7).Value = 0
$ws.Cells.Item(168, 8).Value = 0

$ws.Cells.Item(169, 1).Value = 'Angola'
$ws.Cells.Item(169, 2).Value = 2
$ws.Cells.Item(169, 3).Value = 1
$ws.Cells.Item(169, 4).Value = 0
$ws.Cells.Item(169, 5).Value = 2
$ws.Cells.Item(169, 6).Value = 0
$ws.Cells.Item(169, 7).Value = 0
$ws.Cells.Item(169, 8).Value = 0

$ws.Cells.Item(170, 1).Value = 'Butan'
$ws.Cells.Item(170, 2).Value = 2
$ws.Cells.Item(170, 3).Value = 0
$ws.Cells.Item(170, 4).Value = 0
$ws.Cells.Item(170, 5).Value = 2
$ws.Cells.Item(170, 6).Value = 0
$ws.Cells.Item(170, 7).Value = 0
$ws.Cells.Item(170, 8).Value = 0

$ws.Cells.Item(171, 1).Value = 'Zambia'
$ws.Cells.Item(171, 2).Value = 2
$ws.Cells.Item(171, 3).Value = 0
$ws.Cells.Item(171, 4).Value = 0
$ws.Cells.Item(171, 5).Value = 2
$ws.Cells.Item(171, 6).Value = 0
$ws.Cells.Item(171, 7).Value = 0
$ws.Cells.Item(171, 8).Value = 0

$ws.Cells.Item(172, 1).Value = 'Nicaragua'
$ws.Cells.Item(172, 2).Value = 2
$ws.Cells.Item(172, 3).Value = 0
$ws.Cells.Item(172, 4).Value = 0
$ws.Cells.Item(172, 5).Value = 2
$ws.Cells.Item(172, 6).Value = 0
$ws.Cells.Item(172, 7).Value = 0
$ws.Cells.Item(172, 8).Value = 0

$ws.Cells.Item(173, 1).Value = 'Fiyi'
$ws.Cells.Item(173, 2).Value = 2
$ws.Cells.Item(173, 3).Value = 1
$ws.Cells.Item(173, 4).Value = 0
$ws.Cells.Item(173, 5).Value = 2
$ws.Cells.Item(173, 6).Value = 0
$ws.Cells.Item(173, 7).Value = 0
$ws.Cells.Item(173, 8).Value = 0

$ws.Cells.Item(174, 1).Value = 'Santa Lucia'
$ws.Cells.Item(174, 2).Value = 2
$ws.Cells.Item(174, 3).Value = 0
$ws.Cells.Item(174, 4).Value = 0
$ws.Cells.Item(174, 5).Value = 2
$ws.Cells.Item(174, 6).Value = 0
$ws.Cells.Item(174, 7).Value = 0
$ws.Cells.Item(174, 8).Value = 0

$ws.Cells.Item(176, 1).Value = 'Montserrat'
$ws.Cells.Item(176, 2).Value = 1
$ws.Cells.Item(176, 3).Value = 0
$ws.Cells.Item(176, 4).Value = 0
$ws.Cells.Item(176, 5).Value = 1
$ws.Cells.Item(176, 6).Value = 0
$ws.Cells.Item(176, 7).Value = 0
$ws.Cells.Item(176, 8).Value = 0

$ws.Cells.Item(177, 1).Value = 'Republica de Yibuti'
$ws.Cells.Item(177, 2).Value = 1
$ws.Cells.Item(177, 3).Value = 0
$ws.Cells.Item(177, 4).Value = 0
$ws.Cells.Item(177, 5).Value = 1
$ws.Cells.Item(177, 6).Value = 0
$ws.Cells.Item(177, 7).Value = 0
$ws.Cells.Item(177, 8).Value = 0

$ws.Cells.Item(179, 1).Value = 'Papua Nueva Guinea'
$ws.Cells.Item(179, 2).Value = 1
$ws.Cells.Item(179, 3).Value = 0
$ws.Cells.Item(179, 4).Value = 0
$ws.Cells.Item(179, 5).Value = 1
$ws.Cells.Item(179, 6).Value = 0
$ws.Cells.Item(179, 7).Value = 0
$ws.Cells.Item(179, 8).Value = 0

$ws.Cells.Item(180, 1).Value = 'Antigua y Barbuda'
$ws.Cells.Item(180, 2).Value = 1
$ws.Cells.Item(180, 3).Value = 0
$ws.Cells.Item(180, 4).Value = 0
$ws.Cells.Item(180, 5).Value = 1
$ws.Cells.Item(180, 6).Value = 0
$ws.Cells.Item(180, 7).Value = 0
$ws.Cells.Item(180, 8).Value = 0

$ws.Cells.Item(181, 1).Value = 'Somalia'
$ws.Cells.Item(181, 2).Value = 1
$ws.Cells.Item(181, 3).Value = 0
$ws.Cells.Item(181, 4).Value = 0
$ws.Cells.Item(181, 5).Value = 1
$ws.Cells.Item(181, 6).Value = 0
$ws.Cells.Item(181, 7).Value = 0
$ws.Cells.Item(181, 8).Value = 0

$ws.Cells.Item(182, 1).Value = 'Suazilandia'
$ws.Cells.Item(182, 2).Value = 1
$ws.Cells.Item(182, 3).Value = 0
$ws.Cells.Item(182, 4).Value = 0
$ws.Cells.Item(182, 5).Value = 1
$ws.Cells.Item(182, 6).Value = 0
$ws.Cells.Item(182, 7).Value = 0
$ws.Cells.Item(182, 8).Value = 0

$ws.Cells.Item(183, 1).Value = 'Timor Oriental'
$ws.Cells.Item(183, 2).Value = 1
$ws.Cells.Item(183, 3).Value = 1
$ws.Cells.Item(183, 4).Value = 0
$ws.Cells.Item(183, 5).Value = 1
$ws.Cells.Item(183, 6).Value = 0
$ws.Cells.Item(183, 7).Value = 0
$ws.Cells.Item(183, 8).Value = 0

$ws.Cells.Item(184, 1).Value = 'Gambia'
$ws.Cells.Item(184, 2).Value = 1
$ws.Cells.Item(184, 3).Value = 0
$ws.Cells.Item(184, 4).Value = 0
$ws.Cells.Item(184, 5).Value = 1
$ws.Cells.Item(184, 6).Value = 0
$ws.Cells.Item(184, 7).Value = 0
$ws.Cells.Item(184, 8).Value = 0

$ws.Cells.Item(186, 1).Value = 'San Vicente y las Granadinas'
$ws.Cells.Item(186, 2).Value = 1
$ws.Cells.Item(186, 3).Value = 0
$ws.Cells.Item(186, 4).Value = 0
$ws.Cells.Item(186, 5).Value = 1
$ws.Cells.Item(186, 6).Value = 0
$ws.Cells.Item(186, 7).Value = 0
$ws.Cells.Item(186, 8).Value = 0

$ws.Cells.Item(187, 1).Value = 'Republica del Chad'
$ws.Cells.Item(187, 2).Value = 1
$ws.Cells.Item(187, 3).Value = 0
$ws.Cells.Item(187, 4).Value = 0
$ws.Cells.Item(187, 5).Value = 1
$ws.Cells.Item(187, 6).Value = 0
$ws.Cells.Item(187, 7).Value = 0
$ws.Cells.Item(187, 8).Value = 0

$ws.Cells.Item(188, 1).Value = 'Santa Sede'
$ws.Cells.Item(188, 2).Value = 1
$ws.Cells.Item(188, 3).Value = 0
$ws.Cells.Item(188, 4).Value = 0
$ws.Cells.Item(188, 5).Value = 1
$ws.Cells.Item(188, 6).Value = 0
$ws.Cells.Item(188, 7).Value = 0
$ws.Cells.Item(188, 8).Value = 0
